$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.241.85"
$ws.Range("E2").Value = "  +0.50%  "

# Row 3
$ws.Range("D3").Value = "1.858.98"
$ws.Range("E3").Value = "  +0.82%  "

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "'0.7012"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "

# Row 6
$ws.Range("D6").Value = "'237.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.26%  "

# Row 7
$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").Value = "'0.07925"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.26%  "

# Row 9
$ws.Range("D9").Value = "'0.3051"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.78%  "

# Row 10
$ws.Range("D10").Value = "'23.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.01%  "

# Row 11
$ws.Range("D11").Value = "'0.08188"
$ws.Range("D11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "1.856.55"
$ws.Range("E12").Value = "  +0.54%  "

# Row 13
$ws.Range("D13").Value = "'0.7185"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.66%  "

# Row 14
$ws.Range("D14").Value = "'5.174"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.52%  "

# Row 15
$ws.Range("D15").Value = "'89.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.25%  "

# Row 16
$ws.Range("D16").Value = "29.243.90"
$ws.Range("E16").Value = "  +0.69%  "

# Row 17
$ws.Range("D17").Value = "'5.772"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.01%  "

# Row 18
$ws.Range("D18").Value = "'13.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.88%  "

# Row 19
$ws.Range("D19").Value = "'0.000007785"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.77%  "

# Row 20
$ws.Range("D20").Value = "'236.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.87%  "

# Row 21
$ws.Range("D21").Value = "'0.9995"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.04%  "

# Row 22
$ws.Range("D22").Value = "2.106.30"
$ws.Range("E22").Value = "  +1.33%  "

# Row 23
$ws.Range("D23").Value = "'1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.06%  "

# Row 24
$ws.Range("D24").Value = "'7.450"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.50%  "

# Row 25
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'161.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.06%  "

# Row 26
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'8.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.66%  "

# Row 27
$ws.Range("D27").Value = "'0.1461"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.64%  "

# Row 28
$ws.Range("D28").Value = "'18.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.29%  "

# Row 29
$ws.Range("D29").Value = "'2.010"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.09%  "

# Row 30
$ws.Range("D30").Value = "'1.436"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.45%  "

# Row 31
$ws.Range("D31").Value = "'4.426"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.42%  "

# Row 32
$ws.Range("D32").Value = "'1.482"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.34%  "

# Row 33
$ws.Range("D33").Value = "'4.050"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.32%  "

# Row 34
$ws.Range("D34").Value = "'0.05223"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.53%  "

# Row 35
$ws.Range("D35").Value = "'1.170"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.90%  "

# Row 36
$ws.Range("D36").Value = "'0.7080"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.15%  "

# Row 37
$ws.Range("D37").Value = "'0.9991"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.22%  "

# Row 38
$ws.Range("D38").Value = "'2.674"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.06%  "

# Row 39
$ws.Range("D39").Value = "'0.01846"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.09%  "

# Row 40
$ws.Range("D40").Value = "'2.720"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.96%  "

# Row 41
$ws.Range("D41").Value = "'0.9251"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.49%  "

# Row 42
$ws.Range("D42").Value = "1.142.26"
$ws.Range("E42").Value = "  +8.96%  "

# Row 43
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.927"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.64%  "

# Row 44
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.4277"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.06%  "

# Row 45
$ws.Range("E45").Value = "  +1.28%  "

# Row 46
$ws.Range("D46").Value = "'0.9994"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.02%  "

# Row 47
$ws.Range("D47").Value = "'103.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.17%  "

# Row 48
$ws.Range("D48").Value = "'1.798"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.87%  "

# Row 49
$ws.Range("D49").Value = "2.003.08"
$ws.Range("E49").Value = "  +1.05%  "

# Row 50
$ws.Range("D50").Value = "'9.194"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.16%  "

# Row 51
$ws.Range("E51").Value = "  -1.43%  "
